$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the last existing header cell (AC1) onto the
# three new header cells, then overwrite their text/values.
$headerSrc = $ws.Range("AC1")

$headerSrc.Copy($ws.Range("AD1"))
$ws.Range("AD1").Value = "Wins"

$headerSrc.Copy($ws.Range("AE1"))
$ws.Range("AE1").Value = "Losses"

$headerSrc.Copy($ws.Range("AF1"))
$ws.Range("AF1").Value = "Ties"

# Fill in the win/loss/tie record for every data row.
for ($r = 2; $r -le 56; $r++) {
    $ws.Cells.Item($r, 30).Value = 67
    $ws.Cells.Item($r, 31).Value = 95
    $ws.Cells.Item($r, 32).Value = 0
}
